$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '278.81'
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '6.69%'
$ws.Cells.Item(2, 5).ClearFormats()
$ws.Cells.Item(2, 7).NumberFormat = '@'
$ws.Cells.Item(2, 7).Value = '10'
$ws.Cells.Item(2, 7).ClearFormats()
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '27.30'
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '2.26%'
$ws.Cells.Item(3, 5).ClearFormats()
$ws.Cells.Item(3, 7).NumberFormat = '@'
$ws.Cells.Item(3, 7).Value = '10'
$ws.Cells.Item(3, 7).ClearFormats()
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '4.801'
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '1.97%'
$ws.Cells.Item(4, 5).ClearFormats()
$ws.Cells.Item(4, 7).NumberFormat = '@'
$ws.Cells.Item(4, 7).Value = '10'
$ws.Cells.Item(4, 7).ClearFormats()
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '0.06317'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '1.84%'
$ws.Cells.Item(5, 5).ClearFormats()
$ws.Cells.Item(5, 7).NumberFormat = '@'
$ws.Cells.Item(5, 7).Value = '10'
$ws.Cells.Item(5, 7).ClearFormats()
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '6.931'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '2.94%'
$ws.Cells.Item(6, 5).ClearFormats()
$ws.Cells.Item(6, 7).NumberFormat = '@'
$ws.Cells.Item(6, 7).Value = '10'
$ws.Cells.Item(6, 7).ClearFormats()
$ws.Cells.Item(7, 2).Value = 'GateToken'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '3.372'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '6.26%'
$ws.Cells.Item(7, 5).ClearFormats()
$ws.Cells.Item(7, 7).NumberFormat = '@'
$ws.Cells.Item(7, 7).Value = '10'
$ws.Cells.Item(7, 7).ClearFormats()
$ws.Cells.Item(8, 2).Value = 'MXToken'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.8793'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '3.56%'
$ws.Cells.Item(8, 5).ClearFormats()
$ws.Cells.Item(8, 7).NumberFormat = '@'
$ws.Cells.Item(8, 7).Value = '10'
$ws.Cells.Item(8, 7).ClearFormats()
$ws.Cells.Item(9, 2).Value = 'FTXToken'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.9470'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '3.99%'
$ws.Cells.Item(9, 5).ClearFormats()
$ws.Cells.Item(9, 7).NumberFormat = '@'
$ws.Cells.Item(9, 7).Value = '10'
$ws.Cells.Item(9, 7).ClearFormats()
$ws.Cells.Item(10, 2).Value = 'WazirX'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.1471'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '4.66%'
$ws.Cells.Item(10, 5).ClearFormats()
$ws.Cells.Item(10, 7).NumberFormat = '@'
$ws.Cells.Item(10, 7).Value = '10'
$ws.Cells.Item(10, 7).ClearFormats()
$ws.Cells.Item(11, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.05099'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '0.38%'
$ws.Cells.Item(11, 5).ClearFormats()
$ws.Cells.Item(11, 7).NumberFormat = '@'
$ws.Cells.Item(11, 7).Value = '10'
$ws.Cells.Item(11, 7).ClearFormats()
$ws.Cells.Item(12, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.07320'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '3.18%'
$ws.Cells.Item(12, 5).ClearFormats()
$ws.Cells.Item(12, 7).NumberFormat = '@'
$ws.Cells.Item(12, 7).Value = '10'
$ws.Cells.Item(12, 7).ClearFormats()
$ws.Cells.Item(13, 2).Value = 'BitrueCoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.03146'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '1.20%'
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(13, 7).NumberFormat = '@'
$ws.Cells.Item(13, 7).Value = '10'
$ws.Cells.Item(13, 7).ClearFormats()
$ws.Cells.Item(14, 2).Value = 'BitMartToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.09066'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '0.20%'
$ws.Cells.Item(14, 5).ClearFormats()
$ws.Cells.Item(14, 7).NumberFormat = '@'
$ws.Cells.Item(14, 7).Value = '10'
$ws.Cells.Item(14, 7).ClearFormats()
$ws.Cells.Item(15, 2).Value = 'BitForexToken'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.001557'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '1.21%'
$ws.Cells.Item(15, 5).ClearFormats()
$ws.Cells.Item(15, 7).NumberFormat = '@'
$ws.Cells.Item(15, 7).Value = '10'
$ws.Cells.Item(15, 7).ClearFormats()
$ws.Cells.Item(16, 2).Value = 'One'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.0006269'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '1.42%'
$ws.Cells.Item(16, 5).ClearFormats()
$ws.Cells.Item(16, 7).NumberFormat = '@'
$ws.Cells.Item(16, 7).Value = '10'
$ws.Cells.Item(16, 7).ClearFormats()
$ws.Cells.Item(17, 2).Value = 'TigerCash'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.005862'
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '-1.43%'
$ws.Cells.Item(17, 5).ClearFormats()
$ws.Cells.Item(17, 7).NumberFormat = '@'
$ws.Cells.Item(17, 7).Value = '10'
$ws.Cells.Item(17, 7).ClearFormats()
$ws.Cells.Item(18, 2).Value = 'LEO'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '3.441'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '-0.19%'
$ws.Cells.Item(18, 5).ClearFormats()
$ws.Cells.Item(18, 7).NumberFormat = '@'
$ws.Cells.Item(18, 7).Value = '10'
$ws.Cells.Item(18, 7).ClearFormats()
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '2.291'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '4.74%'
$ws.Cells.Item(19, 5).ClearFormats()
$ws.Cells.Item(19, 7).NumberFormat = '@'
$ws.Cells.Item(19, 7).Value = '10'
$ws.Cells.Item(19, 7).ClearFormats()
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '2.36%'
$ws.Cells.Item(20, 5).ClearFormats()
$ws.Cells.Item(20, 7).NumberFormat = '@'
$ws.Cells.Item(20, 7).Value = '10'
$ws.Cells.Item(20, 7).ClearFormats()
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.1293'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '-1.30%'
$ws.Cells.Item(21, 5).ClearFormats()
$ws.Cells.Item(21, 7).NumberFormat = '@'
$ws.Cells.Item(21, 7).Value = '10'
$ws.Cells.Item(21, 7).ClearFormats()
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '3.884'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '-4.82%'
$ws.Cells.Item(22, 5).ClearFormats()
$ws.Cells.Item(22, 7).NumberFormat = '@'
$ws.Cells.Item(22, 7).Value = '10'
$ws.Cells.Item(22, 7).ClearFormats()
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '2.31%'
$ws.Cells.Item(23, 5).ClearFormats()
$ws.Cells.Item(23, 7).NumberFormat = '@'
$ws.Cells.Item(23, 7).Value = '10'
$ws.Cells.Item(23, 7).ClearFormats()
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '0.001184'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '0.30%'
$ws.Cells.Item(24, 5).ClearFormats()
$ws.Cells.Item(24, 7).NumberFormat = '@'
$ws.Cells.Item(24, 7).Value = '10'
$ws.Cells.Item(24, 7).ClearFormats()
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.004284'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '5.47%'
$ws.Cells.Item(25, 5).ClearFormats()
$ws.Cells.Item(25, 7).NumberFormat = '@'
$ws.Cells.Item(25, 7).Value = '10'
$ws.Cells.Item(25, 7).ClearFormats()
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '-0.10%'
$ws.Cells.Item(26, 5).ClearFormats()
$ws.Cells.Item(26, 7).NumberFormat = '@'
$ws.Cells.Item(26, 7).Value = '10'
$ws.Cells.Item(26, 7).ClearFormats()
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.0001691'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '3.09%'
$ws.Cells.Item(27, 5).ClearFormats()
$ws.Cells.Item(27, 7).NumberFormat = '@'
$ws.Cells.Item(27, 7).Value = '10'
$ws.Cells.Item(27, 7).ClearFormats()
$ws.Cells.Item(28, 7).NumberFormat = '@'
$ws.Cells.Item(28, 7).Value = '10'
$ws.Cells.Item(28, 7).ClearFormats()
$ws.Cells.Item(29, 7).NumberFormat = '@'
$ws.Cells.Item(29, 7).Value = '10'
$ws.Cells.Item(29, 7).ClearFormats()
$ws.Cells.Item(30, 7).NumberFormat = '@'
$ws.Cells.Item(30, 7).Value = '10'
$ws.Cells.Item(30, 7).ClearFormats()
$ws.Cells.Item(31, 7).NumberFormat = '@'
$ws.Cells.Item(31, 7).Value = '10'
$ws.Cells.Item(31, 7).ClearFormats()
$ws.Cells.Item(32, 7).NumberFormat = '@'
$ws.Cells.Item(32, 7).Value = '10'
$ws.Cells.Item(32, 7).ClearFormats()
$ws.Cells.Item(33, 7).NumberFormat = '@'
$ws.Cells.Item(33, 7).Value = '10'
$ws.Cells.Item(33, 7).ClearFormats()
$ws.Cells.Item(34, 7).NumberFormat = '@'
$ws.Cells.Item(34, 7).Value = '10'
$ws.Cells.Item(34, 7).ClearFormats()
$ws.Cells.Item(35, 7).NumberFormat = '@'
$ws.Cells.Item(35, 7).Value = '10'
$ws.Cells.Item(35, 7).ClearFormats()
$ws.Cells.Item(36, 7).NumberFormat = '@'
$ws.Cells.Item(36, 7).Value = '10'
$ws.Cells.Item(36, 7).ClearFormats()
$ws.Cells.Item(37, 7).NumberFormat = '@'
$ws.Cells.Item(37, 7).Value = '10'
$ws.Cells.Item(37, 7).ClearFormats()
$ws.Cells.Item(38, 7).NumberFormat = '@'
$ws.Cells.Item(38, 7).Value = '10'
$ws.Cells.Item(38, 7).ClearFormats()
$ws.Cells.Item(39, 7).NumberFormat = '@'
$ws.Cells.Item(39, 7).Value = '10'
$ws.Cells.Item(39, 7).ClearFormats()
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.04075'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '3.39%'
$ws.Cells.Item(40, 5).ClearFormats()
$ws.Cells.Item(40, 7).NumberFormat = '@'
$ws.Cells.Item(40, 7).Value = '10'
$ws.Cells.Item(40, 7).ClearFormats()
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.006576'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '59.04%'
$ws.Cells.Item(41, 5).ClearFormats()
$ws.Cells.Item(41, 7).NumberFormat = '@'
$ws.Cells.Item(41, 7).Value = '10'
$ws.Cells.Item(41, 7).ClearFormats()
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.1159'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '4.05%'
$ws.Cells.Item(42, 5).ClearFormats()
$ws.Cells.Item(42, 7).NumberFormat = '@'
$ws.Cells.Item(42, 7).Value = '10'
$ws.Cells.Item(42, 7).ClearFormats()
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.002201'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '2.95%'
$ws.Cells.Item(43, 5).ClearFormats()
$ws.Cells.Item(43, 7).NumberFormat = '@'
$ws.Cells.Item(43, 7).Value = '10'
$ws.Cells.Item(43, 7).ClearFormats()
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.01304'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '-1.26%'
$ws.Cells.Item(44, 5).ClearFormats()
$ws.Cells.Item(44, 7).NumberFormat = '@'
$ws.Cells.Item(44, 7).Value = '10'
$ws.Cells.Item(44, 7).ClearFormats()
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.00005227'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '1.24%'
$ws.Cells.Item(45, 5).ClearFormats()
$ws.Cells.Item(45, 7).NumberFormat = '@'
$ws.Cells.Item(45, 7).Value = '10'
$ws.Cells.Item(45, 7).ClearFormats()
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '-0.03%'
$ws.Cells.Item(46, 5).ClearFormats()
$ws.Cells.Item(46, 7).NumberFormat = '@'
$ws.Cells.Item(46, 7).Value = '10'
$ws.Cells.Item(46, 7).ClearFormats()
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '857.83%'
$ws.Cells.Item(47, 5).ClearFormats()
$ws.Cells.Item(47, 7).NumberFormat = '@'
$ws.Cells.Item(47, 7).Value = '10'
$ws.Cells.Item(47, 7).ClearFormats()
$ws.Cells.Item(48, 7).NumberFormat = '@'
$ws.Cells.Item(48, 7).Value = '10'
$ws.Cells.Item(48, 7).ClearFormats()
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.00002101'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '-0.03%'
$ws.Cells.Item(49, 5).ClearFormats()
$ws.Cells.Item(49, 7).NumberFormat = '@'
$ws.Cells.Item(49, 7).Value = '10'
$ws.Cells.Item(49, 7).ClearFormats()
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.0002001'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '-0.03%'
$ws.Cells.Item(50, 5).ClearFormats()
$ws.Cells.Item(50, 7).NumberFormat = '@'
$ws.Cells.Item(50, 7).Value = '10'
$ws.Cells.Item(50, 7).ClearFormats()
$ws.Cells.Item(51, 7).NumberFormat = '@'
$ws.Cells.Item(51, 7).Value = '10'
$ws.Cells.Item(51, 7).ClearFormats()
